# edit.ps1
# Applies the tracked changes described by the commit:
#  1. Removes the (unintended) <w:strike/> character formatting that had
#     been applied to the "About us" paragraphs and the "Most data ..."
#     paragraph.
#  2. Splits the "We would want this..." run into two runs (no visible
#     text change, just a run split after "adopted ").
#  3. Drops the trailing "." from the end of the "... for our customers."
#     sentence (now "... for our customers").
#  4. Removes everything from the blank paragraph right after that
#     paragraph through to the end of the "Data Science Project Steps"
#     section (i.e. "Input:" ... "Model Deployment & Maintenance"), which
#     is no longer part of the document.

$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $para = $d.Paragraphs($paraIndex)
    $rng = $para.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# -- 1a. "Co Wheels is a sustainable transport pioneer ..." (strike removed)
$p4 = '<w:p w14:paraId="0C636B4F" w14:textId="77777777" w:rsidR="008400F0" w:rsidRPr="00DD565F" w:rsidRDefault="008400F0" w:rsidP="00992B02">' +
    '<w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:jc w:val="both"/>' +
    '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="00DD565F"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>Co Wheels is a sustainable transport pioneer which aims to move people from owning cars to using shared transport.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 4 $p4

# -- 1b. "We are one of only three national car clubs ..." (strike removed)
$p5 = '<w:p w14:paraId="5D187C21" w14:textId="77777777" w:rsidR="008400F0" w:rsidRPr="00DD565F" w:rsidRDefault="008400F0" w:rsidP="00992B02">' +
    '<w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:jc w:val="both"/>' +
    '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="00DD565F"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>We are one of only three national car clubs operating in the UK, with market leader Zipcar, which is focused on London, and Enterprise Car Club which also uses its extensive rental office network.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 5 $p5

# -- 1c. "Co Wheels was set up more than 13 years ago ..." (strike removed)
$p6 = '<w:p w14:paraId="6385B218" w14:textId="0F5EF44E" w:rsidR="008400F0" w:rsidRPr="00DD565F" w:rsidRDefault="008400F0" w:rsidP="00992B02">' +
    '<w:pPr><w:jc w:val="both"/>' +
    '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="00DD565F"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>Co Wheels was set up more than 13 years ago in Durham, and now has more than 600 vehicles and 30,000 members based across the UK from Orkney to the Isle of Wight, including 3 cars in Lancaster, and a third of the fleet are EVs. There are two sides to the business, B2C for public cars and B2B closed fleets for employees of bodies like councils, NHS trusts, Universities, and housing providers.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 6 $p6

# -- 2. "We would want this to give us a clearer idea..." split into two runs
$p14 = '<w:p w14:paraId="413618DF" w14:textId="77777777" w:rsidR="008400F0" w:rsidRDefault="008400F0" w:rsidP="00992B02">' +
    '<w:pPr><w:spacing w:before="240" w:after="240"/><w:jc w:val="both"/>' +
    '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t xml:space="preserve">We would want this to give us a clearer idea of profitability in locations and the potential outcomes if we adopted </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>different pricing options and what the impact of this might have on utilisation.</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 14 $p14

# -- 3. "Most data are available ... surveys for our customers" (strike
#       removed from every run + pPr mark, trailing period dropped)
$rsquo = [char]0x2019
$p15 = '<w:p w14:paraId="7D025131" w14:textId="553D8074" w:rsidR="008400F0" w:rsidRPr="0058673B" w:rsidRDefault="008400F0" w:rsidP="00992B02">' +
    '<w:pPr><w:pBdr><w:bottom w:val="single" w:sz="6" w:space="21" w:color="auto"/></w:pBdr><w:jc w:val="both"/>' +
    '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '</w:pPr>' +
    '<w:r w:rsidRPr="0058673B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t xml:space="preserve">Most data </w:t></w:r>' +
    '<w:r w:rsidR="008B66EF" w:rsidRPr="0058673B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>are</w:t></w:r>' +
    '<w:r w:rsidRPr="0058673B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t xml:space="preserve"> available in our booking system TripIQ and can be exported as CSV files or potentially direct from our SQL servers. We also have additional socio demographic data on locations in our GIS system QGIS. It may also need further desk research on competitor data to cost a sample of alternative travel options, as well as studying current </w:t></w:r>' +
    '<w:r w:rsidR="00A40FE6" w:rsidRPr="0058673B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>user' + $rsquo + 's</w:t></w:r>' +
    '<w:r w:rsidRPr="0058673B"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t xml:space="preserve"> surveys or even forming new surveys for our customers</w:t></w:r>' +
    '</w:p>'
Set-ParagraphXml 15 $p15

# -- 4. Remove every paragraph from the blank line after the survey
#       paragraph through to "Model Deployment & Maintenance" (the whole
#       Input / Vital URL / Data Science Project Steps section is gone).
$paraCount = $d.Paragraphs.Count
$startRemove = $d.Paragraphs(16).Range.Start
$endRemove = $d.Paragraphs($paraCount).Range.End
$removeRange = $d.Range($startRemove, $endRemove)
$removeRange.Delete()
